# This edit reorders the weekly price-report data rows (rows 2-18, columns A-T)
# of the active sheet. Row 12 and row 15 keep their original content; all other
# rows are permuted as described by the source diff (rows were re-sequenced,
# most likely due to a re-sort upstream before the data was re-exported).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination (after) row -> source (before) row.
# Row 12 and 15 are identity (unchanged) and included for completeness.
$rowMap = @{
    2  = 10
    3  = 8
    4  = 9
    5  = 17
    6  = 18
    7  = 4
    8  = 11
    9  = 7
    10 = 5
    11 = 16
    12 = 12
    13 = 2
    14 = 3
    15 = 15
    16 = 13
    17 = 14
    18 = 6
}

$firstCol = 1   # A
$lastCol  = 20  # T

# 1) Snapshot every source cell's value (and number vs text nature) up front,
#    since several rows both read-from and write-to overlapping row numbers.
$snapshot = @{}
for ($r = 2; $r -le 18; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $snapshot["$r,$c"] = $cell.Value2
    }
}

# 2) Write the snapshot back out into the destination rows per the mapping.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $val = $snapshot["$srcRow,$c"]
        $ws.Cells.Item($destRow, $c).Value = $val
    }
}
